# The commit inserts a new daily price record for Mango at row 57 of the
# single data sheet, shifting every subsequent record down by one row
# (old row 57 -> new row 58, ..., old row 150 -> new row 151).
#
# Use a native row insert so Excel re-indexes / shifts all the existing
# rows (and the sheet dimension) automatically, then populate the newly
# inserted row with the new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("57:57").Insert()

$ws.Range("A57").Value = 7
$ws.Range("B57").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C57").Value = "Ñuble"
$ws.Range("D57").Value = 45100
$ws.Range("E57").Value = 16
$ws.Range("F57").Value = "Fruta"
$ws.Range("G57").Value = 100108
$ws.Range("H57").Value = "Tropicales y subtropicales"
$ws.Range("I57").Value = 100108002
$ws.Range("J57").Value = "Mango"
$ws.Range("K57").Value = "Sin especificar"
$ws.Range("L57").Value = "Primera"
$ws.Range("M57").Value = 80
$ws.Range("N57").Value = 9000
$ws.Range("O57").Value = 10000
$ws.Range("P57").Value = 9500
$ws.Range("Q57").Value = "$/bandeja 4 kilos"
$ws.Range("R57").Value = "Perú"
$ws.Range("S57").Value = 2375
$ws.Range("T57").Value = 4
